# New version of the "sex" vocabulary: header labels are renamed to the
# new naming convention (label_1_livello_it / label_1_livello_en) and the
# active selection is moved to C7, matching the author's saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 holds the column headers: A1=codice_1 _Llivello, B1=label_ITA_1_livello,
# C1=label_ENG_1_livello. Rename the Italian/English-suffixed header labels to
# the new it/en naming scheme.
$ws.Range("B1").Value = "label_1_livello_it"
$ws.Range("C1").Value = "label_1_livello_en"

# Move/save the active selection to C7 (was B10).
$null = $ws.Range("C7").Select()
